# Fix auto-increase enrollment flags and add debug logging for participation tracking
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New projected values per row (columns C..Q, excluding K which stays 0)
$updates = @{
    2 = @{ C=9456;  D=8390; E=0.8872673434856176; F=0.8852078497573328; G=0.09666309586593051; H=0.08556693124236725; I=41213626.84427914; J=14411981.36879557; L=14411981.36879557; M=55625608.21307472;  N=800758175.3072001; O=783058368.3032;    P=0.01799791973808648; Q=0.01840473450277369 }
    3 = @{ C=9643;  D=8553; E=0.886964637560925;  F=0.8848541278708877; G=0.09541183041090721; H=0.08442555198680833; I=43058552.09083918; J=15064526.98515879; L=15064526.98515879; M=58123079.07599795;  N=837166427.141628;  O=819686251.117558;  P=0.0179946621086971;  Q=0.01837840632878721 }
    4 = @{ C=9833;  D=8760; E=0.8908776568697244; F=0.888618381010347;  G=0.09408327412763762; H=0.08360412673545399; I=45171746.23605794; J=15801201.92710435; L=15801201.92710435; M=60972948.16316229;  N=875679138.500765;  O=858230190.4948111; P=0.0180445110913083;  Q=0.01841137972318848 }
    5 = @{ C=10029; D=8920; E=0.8894206800279191; F=0.8872090710165108; G=0.09311928308151099; H=0.08261627263647085; I=47228158.87806591; J=16496334.53991937; L=16496334.53991937; M=63724493.41798528;  N=914254174.7405434; O=896768069.2775702; P=0.01804348833802248; Q=0.01839531881772808 }
    6 = @{ C=10232; D=9105; E=0.8898553557466771; F=0.8879461673493271; G=0.09199977791416801; H=0.08169085019587474; I=49468210.12644157; J=17245315.50428019; L=17245315.50428019; M=66713525.63072176;  N=955710133.1542411; O=938118307.2810595; P=0.01804450419225281; Q=0.01838287918531528 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

Write-Host "Applied auto-increase enrollment participation updates for rows 2-6"
